$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "59.844.78"
$ws.Range("E2").Value = "  -0.07%  "
$ws.Range("D3").Value = "2.384.08"
$ws.Range("E3").Value = "  -1.35%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'556.52"
$ws.Range("E5").Value = "  +0.90%  "
$ws.Range("D6").Value = "'133.53"
$ws.Range("E6").Value = "  -2.60%  "
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = "  -0.05%  "
$ws.Range("D8").Value = "'0.587"
$ws.Range("E8").Value = "  -0.62%  "
$ws.Range("E9").Value = "  +0.26%  "
$ws.Range("E10").Value = "  -1.65%  "
$ws.Range("E11").Value = "  +1.10%  "
$ws.Range("E12").Value = "  -2.85%  "
$ws.Range("D13").Value = "'24.48"
$ws.Range("E13").Value = "  -4.16%  "
$ws.Range("D14").Value = "2.807.35"
$ws.Range("E14").Value = "  -1.33%  "
$ws.Range("D15").Value = "59.767.32"
$ws.Range("E15").Value = "  -0.09%  "
$ws.Range("E16").Value = "  -0.29%  "
$ws.Range("D17").Value = "2.380.86"
$ws.Range("E17").Value = "  -0.80%  "
$ws.Range("D18").Value = "'11.13"
$ws.Range("E18").Value = "  -1.58%  "
$ws.Range("D19").Value = "'4.49"
$ws.Range("E19").Value = "  +1.80%  "
$ws.Range("D20").Value = "'321.54"
$ws.Range("E20").Value = "  -2.23%  "
$ws.Range("D21").Value = "'6.71"
$ws.Range("E21").Value = "  +0.92%  "
$ws.Range("E22").Value = "  +0.07%  "
$ws.Range("D23").Value = "'64.18"
$ws.Range("E23").Value = "  -3.36%  "
$ws.Range("E24").Value = "  +0.77%  "
$ws.Range("E25").Value = "  +0.01%  "
$ws.Range("D26").Value = "'8.41"
$ws.Range("E26").Value = "  -2.63%  "
$ws.Range("E27").Value = "  -0.08%  "
$ws.Range("E28").Value = "  +1.90%  "
$ws.Range("D29").Value = "0.0₃0760"
$ws.Range("E29").Value = "  -1.82%  "
$ws.Range("D30").Value = "'169.81"
$ws.Range("E30").Value = "  +0.97%  "
$ws.Range("D31").Value = "'6.07"
$ws.Range("E31").Value = "  -0.45%  "
$ws.Range("D32").Value = "'1.10"
$ws.Range("E32").Value = "  +8.77%  "
$ws.Range("D33").Value = "'0.397"
$ws.Range("E33").Value = "  -2.85%  "
$ws.Range("D34").Value = "'18.21"
$ws.Range("E34").Value = "  -2.15%  "
$ws.Range("E36").Value = "  +1.60%  "
$ws.Range("E37").Value = "  +0.02%  "
$ws.Range("D38").Value = "'4.14"
$ws.Range("E38").Value = "  -1.63%  "
$ws.Range("D39").Value = "'318.94"
$ws.Range("E39").Value = "  +1.62%  "
$ws.Range("D40").Value = "'1.58"
$ws.Range("E40").Value = "  -1.48%  "
$ws.Range("D41").Value = "'38.60"
$ws.Range("E41").Value = "  -2.34%  "
$ws.Range("D42").Value = "'145.76"
$ws.Range("E42").Value = "  +5.01%  "
$ws.Range("D43").Value = "'3.53"
$ws.Range("E43").Value = "  -4.18%  "
$ws.Range("E44").Value = "  +0.12%  "
$ws.Range("D45").Value = "'19.70"
$ws.Range("E45").Value = "  +0.97%  "
$ws.Range("E46").Value = "  -1.13%  "
$ws.Range("D47").Value = "'0.572"
$ws.Range("E47").Value = "  -1.52%  "
$ws.Range("D48").Value = "'0.0218"
$ws.Range("E48").Value = "  -2.43%  "
$ws.Range("D49").Value = "'11.06"
$ws.Range("E49").Value = "  +0.09%  "
$ws.Range("E50").Value = "  -0.15%  "
$ws.Range("E51").Value = "  -2.26%  "
